$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 284, shifting existing rows 284..334 down to 285..335
$ws.Rows.Item(284).Insert()

# Populate the new row 284 with the new data record
$ws.Cells.Item(284, 1).Value = 7
$ws.Cells.Item(284, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(284, 3).Value = "Ñuble"
$ws.Cells.Item(284, 4).Value = 45218
$ws.Cells.Item(284, 5).Value = 16
$ws.Cells.Item(284, 6).Value = 100112040
$ws.Cells.Item(284, 7).Value = "Cilantro"
$ws.Cells.Item(284, 8).Value = "Sin especificar"
$ws.Cells.Item(284, 9).Value = "Primera"
$ws.Cells.Item(284, 10).Value = 250
$ws.Cells.Item(284, 11).Value = 2000
$ws.Cells.Item(284, 12).Value = 2000
$ws.Cells.Item(284, 13).Value = 2000
$ws.Cells.Item(284, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(284, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(284, 16).Value = 2000
$ws.Cells.Item(284, 17).Value = 1
$ws.Cells.Item(284, 18).Value = "Hortaliza"
